# Auto-generated script applying the Ragnarok_Profits value updates
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1773
$ws.Range("I2").Value = 1576.6
$ws.Range("J2").Value = 1913.2858
$ws.Range("K2").Value = 1576.6
$ws.Range("L2").Value = 1913.2858
$ws.Range("M2").Value = -1463.6
$ws.Range("N2").Value = -2139.2858
$ws.Range("H3").Value = 42042.6
$ws.Range("J3").Value = 42042.6
$ws.Range("L3").Value = 42042.6
$ws.Range("N3").Value = -42270.6
$ws.Range("H41").Value = 809.35297
$ws.Range("I41").Value = 1255.5
$ws.Range("J41").Value = 412.77777
$ws.Range("K41").Value = 1255.5
$ws.Range("L41").Value = 412.77777
$ws.Range("M41").Value = -815.5
$ws.Range("N41").Value = -1292.77777
$ws.Range("H62").Value = 3948.2144
$ws.Range("I62").Value = 1625.091
$ws.Range("J62").Value = 12466.333
$ws.Range("K62").Value = 1625.091
$ws.Range("L62").Value = 12466.333
$ws.Range("M62").Value = -1001.091
$ws.Range("N62").Value = -13714.333
$ws.Range("H65").Value = 3948.2144
$ws.Range("I65").Value = 1625.091
$ws.Range("J65").Value = 12466.333
$ws.Range("K65").Value = 8125.455
$ws.Range("L65").Value = 62331.665
$ws.Range("M65").Value = -5005.455
$ws.Range("N65").Value = -68571.66500000001
$ws.Range("H76").Value = 9923.857
$ws.Range("I76").Value = 14672.5
$ws.Range("K76").Value = 14672.5
$ws.Range("M76").Value = -14357.5
$ws.Range("H79").Value = 9923.857
$ws.Range("I79").Value = 14672.5
$ws.Range("K79").Value = 14672.5
$ws.Range("M79").Value = -13580.5
$ws.Range("H102").Value = 42042.6
$ws.Range("J102").Value = 42042.6
$ws.Range("L102").Value = 42042.6
$ws.Range("N102").Value = -48532.6
$ws.Range("H112").Value = 4149.7
$ws.Range("J112").Value = 4871.625
$ws.Range("L112").Value = 14614.875
$ws.Range("N112").Value = -16830.875
$ws.Range("H135").Value = 1595.8276
$ws.Range("I135").Value = 747.28
$ws.Range("K135").Value = 6725.52
$ws.Range("M135").Value = -4190.52
$ws.Range("H137").Value = 482699.94
$ws.Range("I137").Value = 1117.5294
$ws.Range("J137").Value = 1301390
$ws.Range("K137").Value = 3352.5882
$ws.Range("L137").Value = 3904170
$ws.Range("M137").Value = -802.5881999999997
$ws.Range("N137").Value = -3909270
$ws.Range("H138").Value = 3267.1072
$ws.Range("I138").Value = 1538.3636
$ws.Range("J138").Value = 4385.706
$ws.Range("K138").Value = 4615.0908
$ws.Range("L138").Value = 13157.118
$ws.Range("M138").Value = 524.9092000000001
$ws.Range("N138").Value = -23437.118
$ws.Range("H140").Value = 110000
$ws.Range("J140").Value = 110000
$ws.Range("L140").Value = 110000
$ws.Range("N140").Value = -120360

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 5150
$ws.Range("J17").Value = 300
$ws.Range("L17").Value = 300
$ws.Range("N17").Value = -646
$ws.Range("H32").Value = 4087.1692
$ws.Range("I32").Value = 4281.2
$ws.Range("J32").Value = 1758.8
$ws.Range("K32").Value = 4281.2
$ws.Range("L32").Value = 1758.8
$ws.Range("M32").Value = -3994.2
$ws.Range("N32").Value = -2332.8
$ws.Range("H61").Value = 5128228
$ws.Range("I61").Value = 5720769.5
$ws.Range("J61").Value = 1671734.9
$ws.Range("K61").Value = 5720769.5
$ws.Range("L61").Value = 1671734.9
$ws.Range("M61").Value = -5720557.5
$ws.Range("N61").Value = -1672158.9
$ws.Range("H63").Value = 3264.4
$ws.Range("I63").Value = 3228.2307
$ws.Range("J63").Value = 3499.5
$ws.Range("K63").Value = 3228.2307
$ws.Range("L63").Value = 3499.5
$ws.Range("M63").Value = -2542.2307
$ws.Range("N63").Value = -4871.5
$ws.Range("H66").Value = 3264.4
$ws.Range("I66").Value = 3228.2307
$ws.Range("J66").Value = 3499.5
$ws.Range("K66").Value = 16141.1535
$ws.Range("L66").Value = 17497.5
$ws.Range("M66").Value = -12709.1535
$ws.Range("N66").Value = -24361.5
$ws.Range("H74").Value = 583461.0600000001
$ws.Range("I74").Value = 659635
$ws.Range("K74").Value = 659635
$ws.Range("M74").Value = -658761
$ws.Range("H77").Value = 583461.0600000001
$ws.Range("I77").Value = 659635
$ws.Range("K77").Value = 3298175
$ws.Range("M77").Value = -3293807
$ws.Range("H132").Value = 1856360.6
$ws.Range("I132").Value = 4747.39
$ws.Range("J132").Value = 7696064
$ws.Range("K132").Value = 14242.17
$ws.Range("L132").Value = 23088192
$ws.Range("M132").Value = -11712.17
$ws.Range("N132").Value = -23093252
$ws.Range("H136").Value = 5128228
$ws.Range("I136").Value = 5720769.5
$ws.Range("J136").Value = 1671734.9
$ws.Range("K136").Value = 17162308.5
$ws.Range("L136").Value = 5015204.699999999
$ws.Range("M136").Value = -17159758.5
$ws.Range("N136").Value = -5020304.699999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1882.12
$ws.Range("I20").Value = 1432.6471
$ws.Range("K20").Value = 1432.6471
$ws.Range("M20").Value = -1185.6471
$ws.Range("H86").Value = 2730.3157
$ws.Range("I86").Value = 1727.9
$ws.Range("J86").Value = 3844.111
$ws.Range("K86").Value = 1727.9
$ws.Range("L86").Value = 3844.111
$ws.Range("M86").Value = -604.9000000000001
$ws.Range("N86").Value = -6090.111
$ws.Range("H89").Value = 2730.3157
$ws.Range("I89").Value = 1727.9
$ws.Range("J89").Value = 3844.111
$ws.Range("K89").Value = 8639.5
$ws.Range("L89").Value = 19220.555
$ws.Range("M89").Value = -3023.5
$ws.Range("N89").Value = -30452.555
$ws.Range("H134").Value = 9093323
$ws.Range("I134").Value = 2222
$ws.Range("K134").Value = 6666
$ws.Range("M134").Value = -4131

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 31568812
$ws.Range("I31").Value = 45457490
$ws.Range("J31").Value = 1013720.8
$ws.Range("K31").Value = 45457490
$ws.Range("L31").Value = 1013720.8
$ws.Range("M31").Value = -45457195
$ws.Range("N31").Value = -1014310.8
$ws.Range("H34").Value = 31568812
$ws.Range("I34").Value = 45457490
$ws.Range("J34").Value = 1013720.8
$ws.Range("K34").Value = 45457490
$ws.Range("L34").Value = 1013720.8
$ws.Range("M34").Value = -45457288
$ws.Range("N34").Value = -1014124.8
$ws.Range("H104").Value = 23947.5
$ws.Range("I104").Value = 25900
$ws.Range("J104").Value = 21995
$ws.Range("K104").Value = 25900
$ws.Range("L104").Value = 21995
$ws.Range("M104").Value = -23279
$ws.Range("N104").Value = -27237

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 14738.6
$ws.Range("J60").Value = 18379.5
$ws.Range("L60").Value = 55138.5
$ws.Range("N60").Value = -55640.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 27.764706
$ws.Range("I2").Value = 30.272728
$ws.Range("J2").Value = 23.166666
$ws.Range("K2").Value = 30.272728
$ws.Range("L2").Value = 23.166666
$ws.Range("M2").Value = 82.727272
$ws.Range("N2").Value = -249.166666
$ws.Range("H70").Value = 14023.3
$ws.Range("I70").Value = 12906.667
$ws.Range("J70").Value = 14501.857
$ws.Range("K70").Value = 12906.667
$ws.Range("L70").Value = 14501.857
$ws.Range("M70").Value = -12636.667
$ws.Range("N70").Value = -15041.857
$ws.Range("H73").Value = 14023.3
$ws.Range("I73").Value = 12906.667
$ws.Range("J73").Value = 14501.857
$ws.Range("K73").Value = 12906.667
$ws.Range("L73").Value = 14501.857
$ws.Range("M73").Value = -11970.667
$ws.Range("N73").Value = -16373.857
$ws.Range("H102").Value = 2990.6
$ws.Range("I102").Value = 2818.1667
$ws.Range("K102").Value = 2818.1667
$ws.Range("M102").Value = -1196.1667
$ws.Range("H122").Value = 3584
$ws.Range("I122").Value = 3913.2812
$ws.Range("J122").Value = 949.75
$ws.Range("K122").Value = 11739.8436
$ws.Range("L122").Value = 2849.25
$ws.Range("M122").Value = -9289.8436
$ws.Range("N122").Value = -7749.25
$ws.Range("H132").Value = 17357666
$ws.Range("I132").Value = 2612
$ws.Range("K132").Value = 7836
$ws.Range("M132").Value = -5306

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1762.5143
$ws.Range("I16").Value = 706.36664
$ws.Range("K16").Value = 706.36664
$ws.Range("M16").Value = -536.36664
$ws.Range("H132").Value = 3472.1724
$ws.Range("I132").Value = 2100.6316
$ws.Range("J132").Value = 6078.1
$ws.Range("K132").Value = 6301.8948
$ws.Range("L132").Value = 18234.3
$ws.Range("M132").Value = -3771.8948
$ws.Range("N132").Value = -23294.3
$ws.Range("H136").Value = 3722.5173
$ws.Range("I136").Value = 3206.375
$ws.Range("J136").Value = 6200
$ws.Range("K136").Value = 9619.125
$ws.Range("L136").Value = 18600
$ws.Range("M136").Value = -7069.125
$ws.Range("N136").Value = -23700

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 99999
$ws.Range("J75").Value = 99999
$ws.Range("L75").Value = 99999
$ws.Range("N75").Value = -101871
$ws.Range("H78").Value = 99999
$ws.Range("J78").Value = 99999
$ws.Range("L78").Value = 299997
$ws.Range("N78").Value = -309357
$ws.Range("H132").Value = 218416.02
$ws.Range("I132").Value = 1001
$ws.Range("K132").Value = 3003
$ws.Range("M132").Value = -473
$ws.Range("H136").Value = 212583.36
$ws.Range("I136").Value = 4333.3555
$ws.Range("J136").Value = 3336333.2
$ws.Range("K136").Value = 13000.0665
$ws.Range("L136").Value = 10008999.6
$ws.Range("M136").Value = -10450.0665
$ws.Range("N136").Value = -10014099.6
